# Adds a new "2022-Q3" sheet (with fund holding detail) right after the
# "总计" (totals) sheet, and inserts a corresponding summary row at the
# top of the "总计" table - matching the commit "feat: add 2022-Q3 data".

function Set-TextCell($ws, $row, $col, $val) {
    # Force the value to be stored as text (keeps leading zeros / trailing
    # decimal zeros exactly as authored, e.g. fund code "003857" or
    # "2.45") the same way typing an apostrophe-prefixed value in Excel
    # keeps it literal instead of letting AutoCorrect coerce it to a number.
    $ws.Cells.Item($row, $col).Value = "'" + $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row for 2022-Q3 right under the header,
#    shifting the existing quarters down, and re-number the running index
#    in column A.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)
$totals.Rows.Item(2).Insert()

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q3"
$totals.Cells.Item(2, 3).Value = 15
$totals.Cells.Item(2, 4).Value = 0.41

for ($r = 3; $r -le 7; $r++) {
    $totals.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q3" worksheet right after "总计" and before the
#    existing "2022-Q2" sheet, then fill in the fund holding detail.
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q3rows = @(
    ,@("003857","前海开源周期优选灵活配置混合A","2.45","93.37","4.82","0.1181",9)
    ,@("001901","前海开源沪港深隆鑫灵活配置混合A","4.36","38.21","1.64","0.0715",9)
    ,@("000969","前海开源大安全核心精选灵活配置混合","1.17","91.57","3.90","0.0456",9)
    ,@("003858","前海开源周期优选灵活配置混合C","0.91","93.37","4.82","0.0439",9)
    ,@("004315","前海开源沪港深新硬件主题灵活配置混合C","1.04","92.20","3.48","0.0362",9)
    ,@("004314","前海开源沪港深新硬件主题灵活配置混合A","0.80","92.20","3.48","0.0278",9)
    ,@("010447","中邮未来成长混合A","0.43","91.79","4.81","0.0207",5)
    ,@("001162","前海开源优势蓝筹股票A","0.48","91.35","3.45","0.0166",8)
    ,@("005328","前海开源价值策略股票","0.34","89.93","3.79","0.0129",8)
    ,@("004320","前海开源沪港深乐享生活灵活配置混合","0.27","71.16","4.48","0.0121",5)
    ,@("014433","国泰智享科技1个月滚动持有混合A","0.10","67.22","3.14","0.0031",10)
    ,@("010448","中邮未来成长混合C","0.06","91.79","4.81","0.0029",5)
    ,@("001638","前海开源优势蓝筹股票C","0.08","91.35","3.45","0.0028",8)
    ,@("001902","前海开源沪港深隆鑫灵活配置混合C","0.04","38.21","1.64","0.0007",9)
    ,@("014434","国泰智享科技1个月滚动持有混合C","0.00","67.22","3.14",0,10)
)

$r = 2
foreach ($row in $q3rows) {
    $q3.Cells.Item($r, 1).Value = $r - 2
    Set-TextCell $q3 $r 2 $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    Set-TextCell $q3 $r 4 $row[2]
    Set-TextCell $q3 $r 5 $row[3]
    Set-TextCell $q3 $r 6 $row[4]
    if ($row[5] -eq 0) {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        Set-TextCell $q3 $r 7 $row[5]
    }
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Restore "总计" as the active sheet (it was active before the edit).
$totals.Activate()

Write-Host "2022-Q3 sheet added"
